$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "T"
$vals = @(0,0.5,1,1.5,2,2.5,3,5,7,7.5,8,8.5,9,9.5,10,10.5,12,12.5,13,13.5,14,14.5,15,15.5,16,18,18.5,18.7,19,19.5,20,20.5,21,24,24.1,24.2,24.28,24.3,25,25.1,25.2)
for ($i=0; $i -lt $vals.Length; $i++) {
  $ws.Columns.Item($i+1).ColumnWidth = $vals[$i]
}
Write-Output "ok"
